$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: convert to numeric 5 / 5 / 5 (C1 was text "6,6")
$ws.Range("A1").Value = 5
$ws.Range("B1").Value = 5
$ws.Range("C1").Value = 5

# Row 2: convert to numeric -5 / -5 / -5 (C2 was text "5,5")
$ws.Range("A2").Value = -5
$ws.Range("B2").Value = -5
$ws.Range("C2").Value = -5

# Row 3: A3 becomes a date-like large int w/ wrap text style, B3/C3 become big ints
$ws.Range("A3").WrapText = $true
$ws.Range("A3").Value = 42767
$ws.Range("B3").Value = 3147483647
$ws.Range("C3").NumberFormat = "General"
$ws.Range("C3").Value = 3147483647

# Row 4: convert to numeric negatives (C4 was text "3,3")
$ws.Range("A4").Value = -42767
$ws.Range("B4").Value = -3147483647
$ws.Range("C4").Value = -3147483647

# Row 5: new row of 5.5 values
$ws.Range("A5").Value = 5.5
$ws.Range("B5").Value = 5.5
$ws.Range("C5").Value = 5.5

# Column widths
$ws.Range("A1:B1").ColumnWidth = 21.49
$ws.Range("C1").ColumnWidth = 27.45

# Selection moves to C4
$ws.Range("C4").Select()
